$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.101.96'
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").Value = '1.907.45'
$ws.Range("E3").Value = '  -1.41%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").Value = '0.7427'
$ws.Range("E5").Value = '  -1.05%  '

$ws.Range("D6").Value = '243.90'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '0.3086'
$ws.Range("E8").Value = '  -3.19%  '

$ws.Range("D9").Value = '26.50'
$ws.Range("E9").Value = '  -5.88%  '

$ws.Range("D10").Value = '0.06981'
$ws.Range("E10").Value = '  -0.93%  '

$ws.Range("D11").Value = '0.08078'

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.962.89'
$ws.Range("E12").Value = '  +1.86%  '

$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.7678'
$ws.Range("E13").Value = '  -2.02%  '

$ws.Range("D14").Value = '5.318'
$ws.Range("E14").Value = '  -1.70%  '

$ws.Range("D15").Value = '92.28'
$ws.Range("E15").Value = '  -0.95%  '

$ws.Range("D16").Value = '14.25'
$ws.Range("E16").Value = '  -1.92%  '

$ws.Range("D17").Value = '30.105.49'
$ws.Range("E17").Value = '  -0.73%  '

$ws.Range("D18").Value = '6.075'
$ws.Range("E18").Value = '  -0.31%  '

$ws.Range("D19").Value = '0.000007825'
$ws.Range("E19").Value = '  -2.24%  '

$ws.Range("D20").Value = '240.04'
$ws.Range("E20").Value = '  -5.05%  '

$ws.Range("D21").Value = '2.186.57'
$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("E23").Value = '  +0.03%  '

$ws.Range("D24").Value = '7.129'
$ws.Range("E24").Value = '  +6.10%  '

$ws.Range("D25").Value = '9.378'
$ws.Range("E25").Value = '  -2.19%  '

$ws.Range("D26").Value = '167.04'
$ws.Range("E26").Value = '  +1.23%  '

$ws.Range("E27").Value = '  -0.61%  '

$ws.Range("E28").Value = '  -2.89%  '

$ws.Range("E29").Value = '  -7.26%  '

$ws.Range("D30").Value = '1.543'
$ws.Range("E30").Value = '  +0.09%  '

$ws.Range("D31").Value = '1.352'
$ws.Range("E31").Value = '  -1.19%  '

$ws.Range("D32").Value = '4.333'
$ws.Range("E32").Value = '  -2.53%  '

$ws.Range("D33").Value = '4.083'
$ws.Range("E33").Value = '  -1.70%  '

$ws.Range("D34").Value = '0.05234'
$ws.Range("E34").Value = '  -0.99%  '

$ws.Range("D35").Value = '1.307'
$ws.Range("E35").Value = '  -2.27%  '

$ws.Range("D36").Value = '0.7492'
$ws.Range("E36").Value = '  -1.25%  '

$ws.Range("E37").Value = '  -2.12%  '

$ws.Range("D38").Value = '0.01963'
$ws.Range("E38").Value = '  -0.12%  '

$ws.Range("D39").Value = '2.798'
$ws.Range("E39").Value = '  -0.11%  '

$ws.Range("D40").Value = '6.335'
$ws.Range("E40").Value = '  -3.19%  '

$ws.Range("D41").Value = '0.4496'
$ws.Range("E41").Value = '  -0.44%  '

$ws.Range("D42").Value = '74.43'

$ws.Range("D43").Value = '1.976'
$ws.Range("E43").Value = '  -0.48%  '

$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.06%  '

$ws.Range("D45").Value = '0.8412'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").Value = '7.741'
$ws.Range("E46").Value = '  +0.19%  '

$ws.Range("D47").Value = '102.00'
$ws.Range("E47").Value = '  +0.14%  '

$ws.Range("D48").Value = '9.897'
$ws.Range("E48").Value = '  -1.28%  '

$ws.Range("D49").Value = '2.073.50'
$ws.Range("E49").Value = '  -1.03%  '

$ws.Range("D50").Value = '36.76'
$ws.Range("E50").Value = '  -2.48%  '

$ws.Range("D51").Value = '0.1181'
$ws.Range("E51").Value = '  -6.38%  '
